$wb = $excel.ActiveWorkbook

$wsProjects = $wb.Worksheets.Item("Projects")
$wsProjects.Activate()

$wsProjects.Range("G3").Value = "https://github.com/brej-30"
$wsProjects.Range("J3").Value = "2025-13"
$wsProjects.Range("A3").Value = "testing-essay-writer"
$wsProjects.Range("B3").Value = "ABC"
$wsProjects.Range("C3").Value = "Generate structured essays with export."
$wsProjects.Range("D3").Value = "A Streamlit app that helps users generate structured essays with consistent sections and export options."
$wsProjects.Range("E3").Value = "Streamlit;LLM;NLP"
$wsProjects.Range("F3").Value = "Python;Streamlit;LangChain"
$wsProjects.Range("H3").Value = "https://essay-writer-app.streamlit.app/"
$wsProjects.Range("I3").Value = "Clean UI workflow;Export options;Prompt tuning"
$wsProjects.Range("K3").Value = 2

$wsProjects.Range("B6").Select()
$wsProjects.Range("A1:K3").SetPhonetic()
